$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New block: rows 2-6 (a third benchmark table, mirroring the
#     existing row 11-14 / row 31-34 blocks) ---------------------------

# Row 2: header numbers 1..5 in B2:F2
$ws.Cells.Item(2, 2).Value = 1
$ws.Cells.Item(2, 3).Value = 2
$ws.Cells.Item(2, 4).Value = 3
$ws.Cells.Item(2, 5).Value = 4
$ws.Cells.Item(2, 6).Value = 5

# Row 3: baseline ("Win64 (low frag)")
$ws.Cells.Item(3, 1).Value = "Win64 (low frag)"
$ws.Cells.Item(3, 6).Value = 515958

# Row 4: "nedmalloc v1.06"
$ws.Cells.Item(4, 1).Value = "nedmalloc v1.06"
$ws.Cells.Item(4, 6).Value = 31089140
$ws.Range("G4").Formula = "=F4/F`$3"

# Row 6: "nedmalloc v1.06 (threadcached sysalloc)" - new string (registered
# first so it lands at shared-string index 7)
$ws.Cells.Item(6, 1).Value = "nedmalloc v1.06 (threadcached sysalloc)"

# Row 5: "nedmalloc v1.06 (patcher)" - new string (registered second, index 8)
$ws.Cells.Item(5, 1).Value = "nedmalloc v1.06 (patcher)"
$ws.Cells.Item(5, 6).Value = 30994083
$ws.Range("G5").Formula = "=F5/F`$3"

# --- Update the existing row 12-14 block with new measurements --------

$ws.Cells.Item(12, 6).Value = 14491780
$ws.Cells.Item(13, 6).Value = 37044111
$ws.Cells.Item(14, 1).Value = "nedmalloc v1.06 (patcher)"
$ws.Cells.Item(14, 6).Value = 36643063

# Row 15 is new: the old "nedmalloc v1.06 (sysalloc)" label moves here
$ws.Cells.Item(15, 1).Value = "nedmalloc v1.06 (sysalloc)"

# --- Selection moves to H13 --------------------------------------------
$ws.Range("H13").Select()
